$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (year 2025) metrics: total_customers, new_customers,
# new_rate and returning_rate. returning_customers (D6) and
# retention_rate (F6) stay the same.
$ws.Range("C6").Value = 437
$ws.Range("E6").Value = 128
$ws.Range("G6").Value = 29.29061784897025
$ws.Range("H6").Value = 70.70938215102976
